$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 01:22"

# 2. Update Brasil row (row 17)
$ws.Range("B17").Value = 25262
$ws.Range("C17").Value = 1832
$ws.Range("E17").Value = 20684
$ws.Range("G17").Value = 204
$ws.Range("H17").Value = 1532

# 3. Update Niger row (row 94)
$ws.Range("B94").Value = 570
$ws.Range("C94").Value = 22
$ws.Range("D94").Value = 90
$ws.Range("E94").Value = 466
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 14

# 4. Update Burkina Faso row (row 95)
$ws.Range("B95").Value = 528
$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 177
$ws.Range("E95").Value = 321
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 30

# 5. Update rows 186-189: San Cristobal y Nieves moves up (now row186),
#    pushing Granada, Curazao, Botsuana each down by one row (rows 187-189).
$ws.Range("A186").Value = "San Cristobal y Nieves"
$ws.Range("B186").Value = 14
$ws.Range("C186").Value = 2
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 14
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

$ws.Range("A187").Value = "Granada"
$ws.Range("B187").Value = 14
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 14
$ws.Range("F187").Value = 2
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

$ws.Range("A188").Value = "Curazao"
$ws.Range("B188").Value = 14
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 10
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 1

$ws.Range("A189").Value = "Botsuana"
$ws.Range("B189").Value = 13
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 12
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 1
